$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "213.14", "1.832.29") are preserved exactly as text and are
# not coerced into floating point numbers / reformatted by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.300.68'
$ws.Range('E2').Value = '  +3.94%  '
$ws.Range('D3').Value = '1.609.63'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '213.14'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').Value = '0.487'
$ws.Range('E7').Value = '  +2.16%  '
$ws.Range('E8').Value = '  +2.54%  '
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('D10').Value = '18.15'
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').Value = '0.0824'
$ws.Range('E11').Value = '  +5.20%  '
$ws.Range('D12').Value = '1.831.92'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').Value = '1.604.25'
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').Value = '0.514'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '26.266.20'
$ws.Range('E16').Value = '  +3.80%  '
$ws.Range('D17').Value = '60.89'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  +2.68%  '
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = '206.69'
$ws.Range('E20').Value = '  +11.28%  '
$ws.Range('D21').Value = '4.27'
$ws.Range('E21').Value = '  +3.04%  '
$ws.Range('D22').Value = '9.37'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '1.77'
$ws.Range('E24').Value = '  +4.97%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '142.37'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('E27').Value = '  -3.80%  '
$ws.Range('D28').Value = '15.30'
$ws.Range('E28').Value = '  +3.03%  '
$ws.Range('D29').Value = '6.48'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('D31').Value = '0.0472'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('D32').Value = '3.16'
$ws.Range('E32').Value = '  +3.21%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +1.36%  '
$ws.Range('E35').Value = '  +2.59%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.0163'
$ws.Range('E36').Value = '  +8.47%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.105.52'
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').Value = '0.782'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = '1.743.74'
$ws.Range('E43').Value = '  +2.40%  '
$ws.Range('D44').Value = '93.15'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0107'
$ws.Range('E46').Value = '  -4.22%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '1.55'
$ws.Range('E47').Value = '  +9.18%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '53.72'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0506'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.409'
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.31%  '
